{"js": "// Add a first-line indent of 720 twips (0.5 in) to the first paragraph\n// of the document body (the \"Hirsh Kabaria\" title line), matching the\n// <w:ind w:firstLine=\"720\"/> inserted into its <w:pPr>.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n// firstLineIndent is expressed in points; 720 twips = 36pt (720 / 20).\nfirstParagraph.firstLineIndent = 36;\n\nawait context.sync();\n", "ps1": "# Add a first-line indent of 720 twips (0.5 in) to the first paragraph\n# of the document body (the \"Hirsh Kabaria\" title line), matching the\n# <w:ind w:firstLine=\"720\"/> inserted into its <w:pPr>.\n$d = $word.ActiveDocument\n$firstParagraph = $d.Paragraphs.Item(1)\n# FirstLineIndent is expressed in points; 720 twips = 36pt (720 / 20).\n$firstParagraph.Format.FirstLineIndent = 36\n"}
